# Workbook: server/LISTAS/mi/BISAGRAS 1842 DISMAY.xlsx (Hoja1)
# Bump the printed date by one day and refresh the price list (per
# "fix bug exeded requeste in google drive").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date in A1 advances one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Updated prices in column D
$ws.Range("D33").Value = 116
$ws.Range("D34").Value = 126.813
$ws.Range("D35").Value = 149.638
$ws.Range("D39").Value = 158.191
$ws.Range("D40").Value = 243.215
